# Add new columns I (I0) and J (IF) with header + data, rows 1-24.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers in I1/J1 - copy the formatting from the existing header cell H1
# (bold font, thin border, centered alignment) so the new headers match
# the rest of the header row, then set their text.
$ws.Cells.Item(1, 8).Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(1, 9).Value = "I0"
$ws.Cells.Item(1, 10).Value = "IF"

# Data rows 2-24: column I (I0) then column J (IF)
$data = @{
    2  = @(7, 7)
    3  = @(7, 8)
    4  = @(7, 9)
    5  = @(2, 4)
    6  = @(9, 9)
    7  = @(11, 12)
    8  = @(6, 6)
    9  = @(6, 7)
    10 = @(7, 9)
    11 = @(7, 8)
    12 = @(1, 7)
    13 = @(1, 3)
    14 = @(1, 5)
    15 = @(1, 4)
    16 = @(1, 7)
    17 = @(1, 5)
    18 = @(1, 5)
    19 = @(1, 6)
    20 = @(1, 6)
    21 = @(4, 7)
    22 = @(2, 7)
    23 = @(8, 9)
    24 = @(5, 6)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
